$wb = $excel.ActiveWorkbook

# --- Add Sheet2 and Sheet3 right after Sheet1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)

# --- Populate Sheet2 ---
$ws2.Range("A1").Value = 0
$ws2.Range("B1").Formula = "= INT(A1/8)"
$ws2.Range("D1").Formula = "=5*24"

$ws2.Range("A2").Formula = "=A1+1"
$ws2.Range("A3").Formula = "=A2+1"

# Fill A4:A42 as one relative-formula fill (creates a shared formula group)
$ws2.Range("A4:A42").Formula = "=A3+1"

# Fill A43:A105 as a second relative-formula fill (second shared formula group)
$ws2.Range("A43:A105").Formula = "=A42+1"

# Fill B2:B65 as one relative-formula fill (creates a shared formula group)
$ws2.Range("B2:B65").Formula = "= INT(A2/8)"

# Fill B66:B105 as a second relative-formula fill (second shared formula group)
$ws2.Range("B66:B105").Formula = "= INT(A66/8)"

# Re-enter A38 individually (it ends up as a standalone formula, not part of the fill)
$ws2.Range("A38").Formula = "=A37+1"

# --- View state ---
$ws2.Activate()
$ws2.Range("F11").Select()

$ws1.Activate()
$ws1.Range("F19").Select()
